# matrici giornalieri con turni da 7 ore
# database aggiornato con il turno di settembre 2016
#
# The "7-hour shift" daily matrices (rows 71-74, OSA_GIORNALIERO block) used
# the 7-hour shift codes "M1"/"P1". They are updated to the plain shift
# codes "M"/"P" (the "R" rest code is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 71 (progressivo 1) ---
$ws.Range("C71").Value = "M"
$ws.Range("D71").Value = "P"
$ws.Range("E71").Value = "M"
$ws.Range("F71").Value = "M"
$ws.Range("G71").Value = "R"
$ws.Range("H71").Value = "P"
$ws.Range("I71").Value = "P"

# --- Row 72 (progressivo 2) ---
$ws.Range("C72").Value = "P"
$ws.Range("D72").Value = "M"
$ws.Range("E72").Value = "P"
$ws.Range("F72").Value = "R"
$ws.Range("G72").Value = "P"
$ws.Range("H72").Value = "M"
$ws.Range("I72").Value = "P"

# --- Row 73 (progressivo 3) ---
$ws.Range("C73").Value = "P"
$ws.Range("D73").Value = "P"
$ws.Range("E73").Value = "R"
$ws.Range("F73").Value = "M"
$ws.Range("G73").Value = "M"
$ws.Range("H73").Value = "P"
$ws.Range("I73").Value = "M"

# --- Row 74 (progressivo 4) ---
$ws.Range("C74").Value = "M"
$ws.Range("D74").Value = "R"
$ws.Range("E74").Value = "P"
$ws.Range("F74").Value = "M"
$ws.Range("G74").Value = "P"
$ws.Range("H74").Value = "M"
$ws.Range("I74").Value = "M"

# Column A widened to fit the updated "OSA_GIORNALIERO" labels.
$ws.Columns.Item(1).ColumnWidth = 17.29

# Leave the selection on the block that was just edited, as the author did.
$ws.Range("C71:I76").Select() | Out-Null
